# Case_2_31 lines_states.xlsx edit:
# Two new contingency lines ("line7", "line8") are inserted into the
# lines/extractions table (pushing the former "extr1".."extr8" rows down by
# two rows), a handful of from_bus/to_bus/in_service values are re-fined,
# and the sheet grows from 15 to 17 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the two brand-new rows (16 & 17) the same column-A formatting
# (bold + border, style index 1) that every other index cell in the table
# uses, by cloning it from the row directly above before filling values in.
$ws.Range("A15").Copy($ws.Range("A16"))
$ws.Range("A15").Copy($ws.Range("A17"))

# name | from_bus | to_bus | in_service
$rows = @(
    @(2,  0, "line1", 7,  9,  $true),
    @(3,  1, "line2", 9,  8,  $true),
    @(4,  2, "line3", 8,  10, $false),
    @(5,  3, "line4", 8,  11, $true),
    @(6,  4, "line5", 10, 5,  $true),
    @(7,  5, "line6", 12, 8,  $true),
    @(8,  6, "line7", 14, 11, $true),
    @(9,  7, "line8", 16, 9,  $true),
    @(10, 8, "extr1", 5,  12, $true),
    @(11, 9, "extr2", 5,  9,  $true),
    @(12, 10, "extr3", 10, 11, $false),
    @(13, 11, "extr4", 7,  8,  $true),
    @(14, 12, "extr5", 9,  11, $true),
    @(15, 13, "extr6", 7,  11, $true),
    @(16, 14, "extr7", 5,  7,  $true),
    @(17, 15, "extr8", 8,  5,  $false)
)

foreach ($r in $rows) {
    $row = $r[0]
    $ws.Range("A$row").Value = $r[1]
    $ws.Range("B$row").Value = $r[2]
    $ws.Range("C$row").Value = $r[3]
    $ws.Range("D$row").Value = $r[4]
    $ws.Range("E$row").Value = $r[5]
}
